# Apply crypto price/volume table updates (commit: "Updated cryptos list on Tue Mar 21 10:34:23 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.025.85'
$ws.Range('E2').Value = '  -0.78%  '
$ws.Range('D3').Value = '1.758.34'
$ws.Range('E3').Value = '  -1.41%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '333.81'
$ws.Range('E5').Value = '  -1.75%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9958'
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3883'
$ws.Range('E7').Value = '  +1.75%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3397'
$ws.Range('E8').Value = '  -1.37%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '45.46'
$ws.Range('E9').Value = '  -3.47%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.124'
$ws.Range('E10').Value = '  -2.50%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07207'
$ws.Range('E11').Value = '  -2.51%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.9983'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '22.33'
$ws.Range('E13').Value = '  -3.82%  '
$ws.Range('E14').Value = '  -4.56%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.039'
$ws.Range('E15').Value = '  -3.60%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.751.91'
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001057'
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06593'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '80.59'
$ws.Range('E19').Value = '  -2.24%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.9960'
$ws.Range('E20').Value = '  -0.44%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.96'
$ws.Range('E21').Value = '  -3.31%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.204'
$ws.Range('E22').Value = '  -3.92%  '
$ws.Range('D23').Value = '28.046.49'
$ws.Range('E23').Value = '  -0.65%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.62'
$ws.Range('E24').Value = '  -3.59%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.382'
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '154.16'
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.92'
$ws.Range('E27').Value = '  -4.06%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.316'
$ws.Range('E28').Value = '  -4.18%  '
$ws.Range('D29').Value = '1.947.89'
$ws.Range('E29').Value = '  -1.68%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.273'
$ws.Range('E30').Value = '  -11.80%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '129.13'
$ws.Range('E31').Value = '  -5.51%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.055'
$ws.Range('E32').Value = '  +3.01%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.840'
$ws.Range('E33').Value = '  -4.86%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.08630'
$ws.Range('E34').Value = '  -2.81%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '12.07'
$ws.Range('E35').Value = '  -5.65%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.142'
$ws.Range('E36').Value = '  -3.22%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06160'
$ws.Range('E37').Value = '  -3.19%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02280'
$ws.Range('E38').Value = '  -6.32%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.6498'
$ws.Range('E39').Value = '  -5.03%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.507'
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.2106'
$ws.Range('E41').Value = '  -3.01%  '
$ws.Range('E42').Value = '  -3.65%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.9960'
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '7.865'
$ws.Range('E44').Value = '  -5.24%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.72'
$ws.Range('E45').Value = '  -2.93%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.817'
$ws.Range('E46').Value = '  -1.70%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5996'
$ws.Range('E47').Value = '  -4.57%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '126.44'
$ws.Range('E48').Value = '  -5.08%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.982'
$ws.Range('E49').Value = '  -5.06%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07001'
$ws.Range('E50').Value = '  -6.10%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.153'
$ws.Range('E51').Value = '  -4.51%  '
